$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Cells.Item(2, 4).Value = '65.784.02'
$ws.Cells.Item(2, 5).Value = '  +1.24%  '

# Row 3: update D3, E3
$ws.Cells.Item(3, 4).Value = '2.700.77'
$ws.Cells.Item(3, 5).Value = '  +2.81%  '

# Row 4: update E4
$ws.Cells.Item(4, 5).Value = '  -0.01%  '

# Row 5: update D5, E5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '608.42'
$ws.Cells.Item(5, 5).Value = '  +2.16%  '

# Row 6: update D6, E6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '157.72'
$ws.Cells.Item(6, 5).Value = '  +1.70%  '

# Row 7: update E7
$ws.Cells.Item(7, 5).Value = '  -0.01%  '

# Row 8: update E8
$ws.Cells.Item(8, 5).Value = '  -0.12%  '

# Row 9: update E9
$ws.Cells.Item(9, 5).Value = '  +5.42%  '

# Row 10: update E10
$ws.Cells.Item(10, 5).Value = '  +3.95%  '

# Row 11: update D11, E11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.402'
$ws.Cells.Item(11, 5).Value = '  +0.37%  '

# Row 12: update E12
$ws.Cells.Item(12, 5).Value = '  +1.24%  '

# Row 13: update E13
$ws.Cells.Item(13, 5).Value = '  +4.66%  '

# Row 14: update D14, E14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.0000203'
$ws.Cells.Item(14, 5).Value = '  +9.75%  '

# Row 15: update D15, E15
$ws.Cells.Item(15, 4).Value = '3.184.84'
$ws.Cells.Item(15, 5).Value = '  +2.78%  '

# Row 16: update D16, E16
$ws.Cells.Item(16, 4).Value = '65.656.00'
$ws.Cells.Item(16, 5).Value = '  +1.21%  '

# Row 17: update D17, E17
$ws.Cells.Item(17, 4).Value = '2.699.61'
$ws.Cells.Item(17, 5).Value = '  +3.39%  '

# Row 18: update D18, E18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '12.71'
$ws.Cells.Item(18, 5).Value = '  +1.70%  '

# Row 19: update E19
$ws.Cells.Item(19, 5).Value = '  +2.36%  '

# Row 20: update D20, E20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '359.10'
$ws.Cells.Item(20, 5).Value = '  +2.22%  '

# Row 21: update E21
$ws.Cells.Item(21, 5).Value = '  +4.09%  '

# Row 22: update D22, E22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.999'
$ws.Cells.Item(22, 5).Value = '  -0.09%  '

# Row 23: update D23, E23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '71.03'
$ws.Cells.Item(23, 5).Value = '  +4.11%  '

# Row 24: update D24, E24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '9.90'
$ws.Cells.Item(24, 5).Value = '  +4.38%  '

# Row 25: update D25, E25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.0000107'
$ws.Cells.Item(25, 5).Value = '  +12.26%  '

# Row 26: update E26
$ws.Cells.Item(26, 5).Value = '  +1.02%  '

# Row 27: update D27, E27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '1.68'
$ws.Cells.Item(27, 5).Value = '  +3.21%  '

# Row 28: update D28, E28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '8.56'
$ws.Cells.Item(28, 5).Value = '  +6.28%  '

# Row 29: update E29
$ws.Cells.Item(29, 5).Value = '  +4.32%  '

# Row 30: update B30, C30, D30, E30
$ws.Cells.Item(30, 2).Value = 'Bittensor'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '550.08'
$ws.Cells.Item(30, 5).Value = '  +5.07%  '

# Row 31: update B31, C31, D31, E31
$ws.Cells.Item(31, 2).Value = 'PancakeSwap'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '2.21'
$ws.Cells.Item(31, 5).Value = '  +5.09%  '

# Row 32: update E32
$ws.Cells.Item(32, 5).Value = '  -0.06%  '

# Row 33: update E33
$ws.Cells.Item(33, 5).Value = '  +3.02%  '

# Row 34: update E34
$ws.Cells.Item(34, 5).Value = '  +6.64%  '

# Row 35: update E35
$ws.Cells.Item(35, 5).Value = '  -1.41%  '

# Row 36: update E36
$ws.Cells.Item(36, 5).Value = '  +2.24%  '

# Row 37: update E37
$ws.Cells.Item(37, 5).Value = '  +3.18%  '

# Row 38: update D38, E38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '164.04'
$ws.Cells.Item(38, 5).Value = '  +0.10%  '

# Row 40: update D40, E40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.00'
$ws.Cells.Item(40, 5).Value = '  +0.10%  '

# Row 41: update D41, E41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '172.10'
$ws.Cells.Item(41, 5).Value = '  +4.46%  '

# Row 42: update E42
$ws.Cells.Item(42, 5).Value = '  -0.04%  '

# Row 43: update D43, E43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '42.64'
$ws.Cells.Item(43, 5).Value = '  +0.93%  '

# Row 44: update E44
$ws.Cells.Item(44, 5).Value = '  +3.05%  '

# Row 45: update E45
$ws.Cells.Item(45, 5).Value = '  +0.78%  '

# Row 46: update D46, E46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '23.67'
$ws.Cells.Item(46, 5).Value = '  +2.90%  '

# Row 47: update E47
$ws.Cells.Item(47, 5).Value = '  +3.00%  '

# Row 48: update E48
$ws.Cells.Item(48, 5).Value = '  +4.97%  '

# Row 49: update D49, E49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.656'
$ws.Cells.Item(49, 5).Value = '  +1.80%  '

# Row 50: update D50, E50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '21.07'
$ws.Cells.Item(50, 5).Value = '  +8.81%  '

# Row 51: update E51
$ws.Cells.Item(51, 5).Value = '  +1.60%  '
